$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Delete the first four data rows (rows 2-5). This shifts the
    # remaining data (old rows 6-20) up into rows 2-16.
    $ws.Range("A2:A5").EntireRow.Delete()

    # Column A is a plain 0-based row index (not a formula), so after
    # the shift it must be renumbered sequentially 0..14 again.
    for ($i = 0; $i -le 14; $i++) {
        $ws.Cells.Item($i + 2, 1).Value = $i
    }
}
